$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.176.15"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "1.900.21"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.59%  "
$ws.Range("D5").Value = "'315.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Value = "'0.5113"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.3937"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "'0.08418"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'42.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").Value = "'1.120"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "'6.245"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "1.899.85"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "'20.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "'7.346"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "'93.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "'0.00001108"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'0.06713"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'17.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'6.015"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("D23").Value = "29.188.26"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'11.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'2.228"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("D26").Value = "2.113.31"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'159.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("D30").Value = "'127.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'1.060"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").Value = "'0.1045"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "'5.915"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "'0.02474"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "'0.06604"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("D37").Value = "'9.073"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("D38").Value = "'0.2190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").Value = "'1.226"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("D40").Value = "'5.106"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").Value = "'0.6460"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "'1.233"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").Value = "'11.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").Value = "'0.6029"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'13.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "'3.672"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "'2.039"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("D49").Value = "'1.229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("D50").Value = "'122.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'1.166"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.28%  "
